$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "updated crypto prices" snapshot (GitHub Actions bot).
# Some Price-column values are plain numeric strings (e.g. "211.85");
# force them through a Text number format so Excel keeps them as text
# (matching the original inlineStr cells) instead of coercing to a
# float, then restore the default "Normal" style so no stray number
# format sticks around on the cell.
$ws.Range("D2").Value = "26.225.56"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.587.51"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.244"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.810.38"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.580.24"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "26.213.21"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "1.421.62"
$ws.Range("E33").Value = "  +8.63%  "
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.587"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.948"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -13.68%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "1.721.37"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -0.14%  "
